$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value2 = 130979935
$ws.Range("B8").Value2 = 91828
$ws.Range("E8").Value2 = 5432
$ws.Range("F8").Value2 = 'Granticka'
$ws.Range("G8").Value2 = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H8").Value2 = ''
$ws.Range("Q8").Value2 = 590645
$ws.Range("R8").Value2 = 6963341
$ws.Range("Z8").Value2 = '10:22'
$ws.Range("AB8").Value2 = '10:22'

# Row 9
$ws.Range("A9").Value2 = 130979926
$ws.Range("B9").Value2 = 79243
$ws.Range("E9").Value2 = 6425
$ws.Range("F9").Value2 = 'Garnlav'
$ws.Range("G9").Value2 = 'Alectoria sarmentosa'
$ws.Range("H9").Value2 = '(Ach.) Ach.'
$ws.Range("Q9").Value2 = 590852
$ws.Range("R9").Value2 = 6963248
$ws.Range("Z9").Value2 = '11:00'
$ws.Range("AB9").Value2 = '11:00'

# Row 18
$ws.Range("A18").Value2 = 130979897
$ws.Range("B18").Value2 = 80348
$ws.Range("E18").Value2 = 6458
$ws.Range("F18").Value2 = 'Lunglav'
$ws.Range("G18").Value2 = 'Lobaria pulmonaria'
$ws.Range("H18").Value2 = '(L.) Hoffm.'
$ws.Range("M18").Value2 = ''
$ws.Range("Q18").Value2 = 590726
$ws.Range("R18").Value2 = 6963153
$ws.Range("Z18").Value2 = '13:24'
$ws.Range("AB18").Value2 = '13:24'
$ws.Range("AC18").Value2 = ''

# Row 19
$ws.Range("A19").Value2 = 130979946
$ws.Range("Q19").Value2 = 590605
$ws.Range("R19").Value2 = 6963364
$ws.Range("Z19").Value2 = '09:47'
$ws.Range("AB19").Value2 = '09:47'

# Row 20
$ws.Range("A20").Value2 = 130979899
$ws.Range("B20").Value2 = 57884
$ws.Range("E20").Value2 = 100109
$ws.Range("F20").Value2 = 'Tretåig hackspett'
$ws.Range("G20").Value2 = 'Picoides tridactylus'
$ws.Range("H20").Value2 = '(Linnaeus, 1758)'
$ws.Range("M20").Value2 = 'färska spår'
$ws.Range("Q20").Value2 = 590850
$ws.Range("R20").Value2 = 6963133
$ws.Range("Z20").Value2 = '13:16'
$ws.Range("AB20").Value2 = '13:16'
$ws.Range("AC20").Value2 = 'färska ringhack på gran'

# Row 21
$ws.Range("A21").Value2 = 130979947
$ws.Range("B21").Value2 = 91808
$ws.Range("E21").Value2 = 1202
$ws.Range("F21").Value2 = 'Ullticka'
$ws.Range("G21").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H21").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q21").Value2 = 590591
$ws.Range("R21").Value2 = 6963354
$ws.Range("Z21").Value2 = '09:45'
$ws.Range("AB21").Value2 = '09:45'
